# Update the "Förändrad" column (column C) for all data rows from the
# old date serial (46060 / 2026-02-07) to the new date serial
# (46061 / 2026-02-08). Data rows run from row 2 to row 369.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 369; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
